$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

$ws.Range("B13").Value = "shuttle"
$ws.Range("B18").Value = "myself"
